$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: was 44533/100/2000/2200/2100/2100 -> now 44608/120/600/650/625/625
$ws.Range("D3").Value = 44608
$ws.Range("J3").Value = 120
$ws.Range("K3").Value = 600
$ws.Range("L3").Value = 650
$ws.Range("M3").Value = 625
$ws.Range("P3").Value = 625

# Row 5: was 44754/200/700/750/725/725 -> now 44624/120/650/700/675/675
$ws.Range("D5").Value = 44624
$ws.Range("J5").Value = 120
$ws.Range("K5").Value = 650
$ws.Range("L5").Value = 700
$ws.Range("M5").Value = 675
$ws.Range("P5").Value = 675

# Row 6: was 44624/120/650/700/675/675 -> now 44533/100/2000/2200/2100/2100
$ws.Range("D6").Value = 44533
$ws.Range("J6").Value = 100
$ws.Range("K6").Value = 2000
$ws.Range("L6").Value = 2200
$ws.Range("M6").Value = 2100
$ws.Range("P6").Value = 2100

# Row 7: was 44608/120/600/650/625/625 -> now 44754/200/700/750/725/725
$ws.Range("D7").Value = 44754
$ws.Range("J7").Value = 200
$ws.Range("K7").Value = 700
$ws.Range("L7").Value = 750
$ws.Range("M7").Value = 725
$ws.Range("P7").Value = 725
